$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.264.62"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "1.447.32"
$ws.Range("E3").Value = "  +2.72%  "

$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").Value = "'0.9292"
$ws.Range("E5").Value = "  -7.16%  "

$ws.Range("D6").Value = "'277.03"
$ws.Range("E6").Value = "  +1.72%  "

$ws.Range("D7").Value = "'0.3669"
$ws.Range("E7").Value = "  -0.76%  "

$ws.Range("D8").Value = "'0.3126"
$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("D9").Value = "'38.89"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").Value = "'1.022"
$ws.Range("E10").Value = "  +3.26%  "

$ws.Range("D11").Value = "'0.06529"
$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").Value = "'0.9997"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "'5.405"
$ws.Range("E13").Value = "  +1.94%  "

$ws.Range("D14").Value = "'17.57"
$ws.Range("E14").Value = "  +4.05%  "

$ws.Range("D15").Value = "'6.095"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").Value = "'0.00001017"
$ws.Range("E17").Value = "  +1.15%  "

$ws.Range("D18").Value = "'0.9399"
$ws.Range("E18").Value = "  -6.09%  "

$ws.Range("D19").Value = "'0.05617"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("D20").Value = "'67.32"
$ws.Range("E20").Value = "  -7.80%  "

$ws.Range("D21").Value = "'5.422"
$ws.Range("E21").Value = "  -2.59%  "

$ws.Range("D22").Value = "'14.45"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").Value = "'2.274"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("D25").Value = "20.265.82"
$ws.Range("E25").Value = "  +1.50%  "

$ws.Range("D26").Value = "'2.191"
$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").Value = "'135.33"
$ws.Range("E27").Value = "  -1.85%  "

$ws.Range("D28").Value = "'16.97"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").Value = "1.601.64"
$ws.Range("E29").Value = "  +2.08%  "

$ws.Range("D30").Value = "'110.38"
$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("D31").Value = "'3.759"
$ws.Range("E31").Value = "  -2.81%  "

$ws.Range("D32").Value = "'0.8141"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "'4.856"
$ws.Range("E33").Value = "  -7.70%  "

$ws.Range("D34").Value = "'0.07657"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").Value = "'1.505"
$ws.Range("E35").Value = "  +17.98%  "

$ws.Range("D36").Value = "'0.05985"
$ws.Range("E36").Value = "  +3.29%  "

$ws.Range("D37").Value = "'4.692"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").Value = "'1.135"
$ws.Range("E38").Value = "  +7.15%  "

$ws.Range("D39").Value = "'10.27"
$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").Value = "'0.01989"
$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("D41").Value = "'0.9412"
$ws.Range("E41").Value = "  -5.97%  "

$ws.Range("D42").Value = "'0.1827"
$ws.Range("E42").Value = "  -5.92%  "

$ws.Range("D43").Value = "'7.094"
$ws.Range("E43").Value = "  -15.73%  "

$ws.Range("D44").Value = "'0.5247"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").Value = "'3.523"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("D46").Value = "'11.90"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").Value = "'120.05"
$ws.Range("E47").Value = "  +8.75%  "

$ws.Range("D48").Value = "'0.5156"
$ws.Range("E48").Value = "  +1.35%  "

$ws.Range("D49").Value = "'1.769"
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("D51").Value = "'0.9918"
$ws.Range("E51").Value = "  -0.82%  "

